$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its current location
#    (the empty paragraph right after "All-in at the River" block).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "Version: ..." paragraph (last real text paragraph,
#    right before the section break) and update it:
#      - add wordWrap="0" to its paragraph formatting
#      - change the text from "Version: 2021/10" to "Version: EDG"
#      - re-insert the "_GoBack" bookmark at the very end of its text
# ------------------------------------------------------------------
$versionPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Version:*") {
        $versionPara = $cand
        break
    }
}

if ($versionPara -eq $null) {
    throw "Could not locate the 'Version:' paragraph"
}

# Turn off word wrap for this paragraph (adds <w:wordWrap w:val="0"/>)
$versionPara.Format.WordWrap = $false

$paraStart = $versionPara.Range.Start
$paraEnd = $versionPara.Range.End

# Replace the "10" run (the characters right before the paragraph mark)
# with "EDG".
$tailRange = $d.Range($paraEnd - 3, $paraEnd - 1)
if ($tailRange.Text -ne "10") {
    throw "Unexpected tail text: [" + $tailRange.Text + "]"
}
$tailRange.Text = "EDG"

# Replace the "Version: 2021/" run with "Version: ".
$headRange = $d.Range($paraStart, $paraEnd - 3)
if ($headRange.Text -ne "Version: 2021/") {
    throw "Unexpected head text: [" + $headRange.Text + "]"
}
$headRange.Text = "Version: "

# Recompute the end of the paragraph's text (just before the
# paragraph mark) after the edits above.
$newParaEnd = $versionPara.Range.End

# Adding a bookmark as a truly zero-length range right before a
# paragraph mark is unreliable, so insert a temporary placeholder
# character, anchor the bookmark around it (non-collapsed range),
# then delete the placeholder -- the bookmark collapses correctly
# to the boundary and survives.
$insertPos = $newParaEnd - 1
$placeholderRange = $d.Range($insertPos, $insertPos)
$placeholderRange.InsertAfter("X")

$wrapRange = $d.Range($insertPos, $insertPos + 1)
$wrapRange.Bookmarks.Add("_GoBack")

$placeholderRange2 = $d.Range($insertPos, $insertPos + 1)
$placeholderRange2.Text = ""
